$wb = $excel.ActiveWorkbook

# --- Sheet "乳品" (first sheet): remove the stray empty row 5 ---
$ws1 = $wb.Worksheets.Item("乳品")
$ws1.Rows.Item(5).Delete()

# --- Sheet "水果類" (fifth sheet): update row 3, add rows 6 and 7 ---
$ws5 = $wb.Worksheets.Item("水果類")

# helper: write a numeric-looking string into a cell while keeping it as TEXT
# (force Text format so Excel doesn't auto-convert to a number, then reset
# the cell's style back to the default/normal style so no stray style index
# is left referenced on the cell)
$fmtSource = $ws5.Range("A2")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $fmtSource.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

Set-TextValue $ws5.Range("B3") "10"
$ws5.Range("C3").Value = "公斤"

$ws5.Range("A6").Value = "香蕉"
Set-TextValue $ws5.Range("B6") "10"
$ws5.Range("C6").Value = "公克"

$ws5.Range("A7").Value = "芒果"
Set-TextValue $ws5.Range("B7") "37"
$ws5.Range("C7").Value = "台斤"
